# feat: add 2022-Q3 data
#
# Before: sheet "总计" (summary) + sheet "2022-Q2" (fund detail for that quarter).
# After : sheet "总计" (summary, now with 2 history rows) + sheet "2022-Q3" (new
#          fund detail, re-using the old sheet's position/rId) + sheet "2022-Q2"
#          (the OLD fund detail, preserved verbatim on a brand-new sheet).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)     # "总计"
$q2    = $wb.Worksheets.Item(2)     # currently "2022-Q2"

# ------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet so the OLD fund data keeps
#    living on, unchanged, right after itself. Worksheet.Copy clones
#    values/styles/sheet options (outline/pageSetup/margins) exactly.
# ------------------------------------------------------------------
$q2.Copy($null, $q2)

# The original sheet (still rId2 / same physical part) becomes the new
# latest-quarter sheet; free up the "2022-Q2" name for the fresh copy.
$q2.Name = "2022-Q3"
$q3 = $q2

$q2Archive = $wb.Worksheets.Item(3)
$q2Archive.Name = "2022-Q2"

# ------------------------------------------------------------------
# 2) Replace the "2022-Q3" sheet's content with the new quarter's fund
#    data (wipes values AND old formatting/margins).
# ------------------------------------------------------------------
$q3.Cells.Clear()

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "011729"
$q3.Range("B2").Style = "Normal"
$q3.Range("C2").Value = "工银聚享混合A"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "1.36"
$q3.Range("D2").Style = "Normal"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "26.62"
$q3.Range("E2").Style = "Normal"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "0.93"
$q3.Range("F2").Style = "Normal"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0126"
$q3.Range("G2").Style = "Normal"
$q3.Range("H2").Value = 8

$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "011730"
$q3.Range("B3").Style = "Normal"
$q3.Range("C3").Value = "工银聚享混合C"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "0.00"
$q3.Range("D3").Style = "Normal"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "26.62"
$q3.Range("E3").Style = "Normal"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "0.93"
$q3.Range("F3").Style = "Normal"
$q3.Range("G3").Value = 0
$q3.Range("H3").Value = 8

# style index 2 (bold/centered/thin-border "label" look) -- reuse the
# same style already used throughout the "总计" sheet so no duplicate
# style gets created.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2:A3").PasteSpecial(-4122)

# page margins for the new sheet follow the "总计" sheet's convention
# (0.75in/1in/0.5in) rather than the old fund-detail sheet's (0.7/0.75/0.3).
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# 3) Update the "总计" summary sheet: relabel the latest row as Q3 with
#    its own figures, and push the old Q2 totals down into a new row.
# ------------------------------------------------------------------
$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.18
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4) Restore "总计" as the active sheet/tab (matches the unchanged
#    bookViews / activeTab="0" in the target workbook).
# ------------------------------------------------------------------
$total.Select()
